# Update the weekly salary distribution numbers for team members
# (rows 8-11, column B) and move the active selection to B12, matching
# the author's latest upload of the tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value  = 110
$ws.Range("B9").Value  = 90
$ws.Range("B10").Value = 90
$ws.Range("B11").Value = 110

# Leave the selection where the author left it when they saved the file.
$ws.Range("B12").Select()
